$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workblocks")

# Delete the two "SuppressSuccessful" rows (rows 7 then 6), from bottom up
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()

# Update remaining rows with the simplified workblock interface
$ws.Range("A2").Value = "wbInitAllApplications_Type"
$ws.Range("B2").Value = "Init"
$ws.Range("C2").Value = "Name of Workblock"

$ws.Range("A3").Value = "wbGetTransactionData_Type"
$ws.Range("B3").Value = "GetData"
$ws.Range("C3").Value = "Name of Workblock"

$ws.Range("A4").Value = "wbProcessTransaction_Type"
$ws.Range("B4").Value = "Process"
$ws.Range("C4").Value = "Name of Workblock"

$ws.Range("A5").Value = "wbCloseAllApplications_Type"
$ws.Range("B5").Value = "Close"
$ws.Range("C5").Value = "Name of Workblock"

# Make the Workblocks sheet the active/selected tab, with B13 selected
$ws.Activate()
$ws.Range("B13").Select()
